$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 42608.88890046296

$ws.Range("B8").Value = -2
$ws.Range("C8").Value = 51
$ws.Range("D8").Value = 45
$ws.Range("E8").Value = 30
$ws.Range("F8").Value = 69
$ws.Range("G8").Value = 27300
$ws.Range("H8").Value = 23374
$ws.Range("I8").Value = 1363
$ws.Range("J8").Value = 270
$ws.Range("K8").Value = 237
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 9
$ws.Range("N8").Value = "Named"
